# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 213
$wsOff.Range("C3").Value = 148
$wsOff.Range("D3").Value = 53
$wsOff.Range("E3").Value = 32
$wsOff.Range("F3").Value = 7

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 187
$wsDef.Range("C3").Value = 116
$wsDef.Range("D3").Value = 40
$wsDef.Range("E3").Value = 16
